# Update test data for login
# The "login.feature" test cases (TC_001 "User verify login with valid
# credentials" and TC_002 "User tries to log in with invalid credentials")
# are now flagged as both a SmokeTest and a RegressionTest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Yes"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = "Yes"

# Move the active cell/selection as recorded in the saved workbook.
$ws.Range("G12").Select()
